# Append more rows of exported NLQ data to Sheet1, continuing the existing
# series (columns: A=datetime serial, B=symbol, C-F=price fields, G=volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 317 already exists but is missing its "symbol" (B) cell - fill it in
# so it matches the rest of the column before using it as the template row
# for the newly appended rows (this also keeps formatting/style consistent).
$ws.Cells.Item(317, 2).Value = "ECONOMICS:CNCBBS"

# New rows 318-326 continue the same series; copy row 317 as a formatting
# template (date style on column A, etc.) into each new row, then
# overwrite with the actual values for that row.
for ($row = 318; $row -le 326; $row++) {
    $ws.Range("A317:G317").Copy($ws.Range("A" + $row + ":G" + $row))
}

$data = @(
    @{ Row = 318; Date = 45230; Val = 43325980000000; Symbol = $true },
    @{ Row = 319; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 320; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 321; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 322; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 323; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 324; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 325; Date = 45257; Val = 44065463000000; Symbol = $true },
    @{ Row = 326; Date = 45257; Val = 44065463000000; Symbol = $false }
)

foreach ($r in $data) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date

    if ($r.Symbol) {
        $ws.Cells.Item($row, 2).Value = "ECONOMICS:CNCBBS"
    } else {
        $ws.Cells.Item($row, 2).ClearContents()
    }

    $ws.Cells.Item($row, 3).Value = $r.Val
    $ws.Cells.Item($row, 4).Value = $r.Val
    $ws.Cells.Item($row, 5).Value = $r.Val
    $ws.Cells.Item($row, 6).Value = $r.Val
    $ws.Cells.Item($row, 7).Value = 0
}
